# Y4_B2526_Excuses.xlsx - attendance app upload
# - Row 2 (A2/C2): student ID + log date were corrected for a different
#   student/log entry.
# - Rows 3-4: the two extra "general surgery" excuse rows were removed,
#   leaving just the header + the single (updated) data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A2 ("Student ID") and C2 ("Log Date") ------------------------
# These columns hold ID/date-shaped *text*, not numbers/dates, in the
# source file. A plain `.Value = "211741"` / `.Value = "07/09/2025"`
# assignment gets "smart" auto-converted by Excel into a real number /
# date serial (and drags in a brand-new number-format style along with
# it). Writing a text-literal formula and collapsing it back to a value
# via copy / paste-special-values keeps the cell's existing style intact
# and stores the literal text, exactly like the original inline string.
$ws.Range("A2").Formula = '="211741"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("C2").Formula = '="07/09/2025"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues

# --- Remove the two trailing rows ----------------------------------------
$ws.Rows("3:4").Delete()
